$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54: Downlink Lost.
$ws.Range("D54").Value = "[0, 0, 0, 0, 0, 1, 0]"
$ws.Range("E54").Value = "['CommunicationIssue']"

# Row 61: Exiting GPS mode: Unknown Error
$ws.Range("D61").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E61").Value = "['SoftwareFault']"

# Row 73: IMU attitude restricted. Ensure aircraft is level
$ws.Range("D73").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'HardwareFault']"
